$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = " "
$ws.Range("A6").Value = " "
$ws.Range("A7").Value = " "
